$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C10 ("RET1" row) value changes from 18 to 1
$ws.Range("C10").Value = 1
